$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$data = @{}
$data[2] = @(2,1,1.3106835,2.621367,0.03638853870331563,0.02494331542607714,2,1,8.7735875,17.547175,0.171126311618442,0.1219656100585016,11.49939637205625,45.997585488225,0.006227036413483327,0.003042226682823132)
$data[3] = @(2,1,1.3106835,2.621367,0.03638853870331563,0.02494331542607714,3,1,1.161330666666667,3.483992,0.02265142207290381,0.02421627468347122,1.522136942844,9.132821657064,0.0008242521487849982,0.0006040341778743489)
$data[4] = @(2,1,1.3106835,2.621367,0.03638853870331563,0.02494331542607714,3,1,7.719567000000001,23.158701,0.1505679436150197,0.1609697911844745,10.1179090940445,60.70745456426701,0.005478947443713788,0.004015120275584118)
$data[5] = @(2,1,1.3106835,2.621367,0.03638853870331563,0.02494331542607714,3,1,16.007757,48.023271,0.3122267160034931,0.3337966108230954,20.9811029719095,125.886617831457,0.01136147393950224,0.008325994151915981)
$data[6] = @(2,1,1.3106835,2.621367,0.03638853870331563,0.02494331542607714,3,1,16.441887,49.325661,0.3206942973278844,0.3428491671966479,21.5501099997645,129.300659998587,0.01166959685024833,0.008551794920953846)
$data[7] = @(2,1,1.3106835,2.621367,0.03638853870331563,0.02494331542607714,2,1,1.165529,2.331058,0.02273330936225701,0.01620254605380927,1.5276396290715,6.110558516286001,0.0008272319075829368,0.0004041452169257059)
$data[8] = @(3,1,1.420479333333333,4.261438,0.03943680316282304,0.04054922191462366,2,1,8.7735875,17.547175,0.171126311618442,0.1219656100585016,12.46269972294167,74.77619833765,0.006748674667276414,0.004945610588214637)
$data[9] = @(3,1,1.420479333333333,4.261438,0.03943680316282304,0.04054922191462366,3,1,1.161330666666667,3.483992,0.02265142207290381,0.02421627468347122,1.649646211166222,14.846815900496,0.0008932996736471327,0.0009819510960855575)
$data[10] = @(3,1,1.420479333333333,4.261438,0.03943680316282304,0.04054922191462366,3,1,7.719567000000001,23.158701,0.1505679436150197,0.1609697911844745,10.965485385782,98.689368472038,0.005937918354976569,0.006527199784289888)
$data[11] = @(3,1,1.420479333333333,4.261438,0.03943680316282304,0.04054922191462366,3,1,16.007757,48.023271,0.3122267160034931,0.3337966108230954,22.738687991522,204.648191923698,0.01231322354120441,0.01353519284661497)
$data[12] = @(3,1,1.420479333333333,4.261438,0.03943680316282304,0.04054922191462366,3,1,16.441887,49.325661,0.3206942973278844,0.3428491671966479,23.355360684502,210.198246160518,0.01264715787915962,0.01390226696390079)
$data[13] = @(3,1,1.420479333333333,4.261438,0.03943680316282304,0.04054922191462366,2,1,1.165529,2.331058,0.02273330936225701,0.01620254605380927,1.655609856900667,9.933659141404,0.000896529046558892,0.000657000635517822)
$data[14] = @(3,1,17.64166233333333,52.924987,0.4897859114022,0.5036016111677731,2,1,8.7735875,17.547175,0.171126311618442,0.1219656100585016,154.7806681269542,928.684008761725,0.08381525650093549,0.06142207773252176)
$data[15] = @(3,1,17.64166233333333,52.924987,0.4897859114022,0.5036016111677731,3,1,1.161330666666667,3.483992,0.02265142207290381,0.02421627468347122,20.48780347867822,184.390231308104,0.0110943474045331,0.01219535494707746)
$data[16] = @(3,1,17.64166233333333,52.924987,0.4897859114022,0.5036016111677731,3,1,7.719567000000001,23.158701,0.1505679436150197,0.1609697911844745,136.185994373543,1225.673949361887,0.07374605749143746,0.08106464618984136)
$data[17] = @(3,1,17.64166233333333,52.924987,0.4897859114022,0.5036016111677731,3,1,16.007757,48.023271,0.3122267160034931,0.3337966108230954,282.403443708053,2541.630993372477,0.1529242466618867,0.168100511012853)
$data[18] = @(3,1,17.64166233333333,52.924987,0.4897859114022,0.5036016111677731,3,1,16.441887,49.325661,0.3206942973278844,0.3428491671966479,290.0622185768229,2610.559967191407,0.157071548698226,0.1726593929877611)
$data[19] = @(3,1,17.64166233333333,52.924987,0.4897859114022,0.5036016111677731,2,1,1.165529,2.331058,0.02273330936225701,0.01620254605380927,20.56186905770767,123.371214346246,0.01113445464518122,0.008159628297718392)
$data[20] = @(3,1,13.46419533333333,40.39258599999999,0.3738067908812475,0.3843509945280257,2,1,8.7735875,17.547175,0.171126311618442,0.1219656100585016,118.1292958740916,708.7757752445499,0.06396817738143414,0.04687760352420247)
$data[21] = @(3,1,13.46419533333333,40.39258599999999,0.3738067908812475,0.3843509945280257,3,1,1.161330666666667,3.483992,0.02265142207290381,0.02421627468347122,15.63638294259022,140.727446483312,0.008467255393968828,0.009307549258356015)
$data[22] = @(3,1,13.46419533333333,40.39258599999999,0.3738067908812475,0.3843509945280257,3,1,7.719567000000001,23.158701,0.1505679436150197,0.1609697911844745,103.937757976754,935.4398217907859,0.05628331981231912,0.0618688993307214)
$data[23] = @(3,1,13.46419533333333,40.39258599999999,0.3738067908812475,0.3843509945280257,3,1,16.007757,48.023271,0.3122267160034931,0.3337966108230954,215.531567096534,1939.784103868806,0.1167124667366564,0.1282950593399411)
$data[24] = @(3,1,13.46419533333333,40.39258599999999,0.3738067908812475,0.3843509945280257,3,1,16.441887,49.325661,0.3206942973278844,0.3428491671966479,221.3767782165939,1992.391003949346,0.1198777061380531,0.131774418385137)
$data[25] = @(3,1,13.46419533333333,40.39258599999999,0.3738067908812475,0.3843509945280257,2,1,1.165529,2.331058,0.02273330936225701,0.01620254605380927,15.69291012266466,94.157460735988,0.008497865418815913,0.006227464689667731)
$data[26] = @(3,1,0.5283693333333334,1.585108,0.01466913098062581,0.01508291239967478,2,1,8.7735875,17.547175,0.171126311618442,0.1219656100585016,4.635694578316667,27.8141674699,0.002510274279362315,0.001839596612285273)
$data[27] = @(3,1,0.5283693333333334,1.585108,0.01466913098062581,0.01508291239967478,3,1,1.161330666666667,3.483992,0.02265142207290381,0.02421627468347122,0.6136115101262222,5.522503591135999,0.0003322766772848647,0.0003652519496972585)
$data[28] = @(3,1,0.5283693333333334,1.585108,0.01466913098062581,0.01508291239967478,3,1,7.719567000000001,23.158701,0.1505679436150197,0.1609697911844745,4.078782469412,36.709042224708,0.002208700886372206,0.00242789325942937)
$data[29] = @(3,1,0.5283693333333334,1.585108,0.01466913098062581,0.01508291239967478,3,1,16.007757,48.023271,0.3122267160034931,0.3337966108230954,8.458007894252,76.122071048268,0.004580094592705897,0.005034625040353081)
$data[30] = @(3,1,0.5283693333333334,1.585108,0.01466913098062581,0.01508291239967478,3,1,16.441887,49.325661,0.3206942973278844,0.3428491671966479,8.687388872931999,78.186499856388,0.004704306652242495,0.005171163955128491)
$data[31] = @(3,1,0.5283693333333334,1.585108,0.01466913098062581,0.01508291239967478,2,1,1.165529,2.331058,0.02273330936225701,0.01620254605380927,0.6158297807106667,3.694978684264,0.0003334778926580352,0.0002443815827813015)
$data[32] = @(2,1,1.65374,3.30748,0.0459128248697883,0.03147194456382552,2,1,8.7735875,17.547175,0.171126311618442,0.1219656100585016,14.50923259225,58.036930369,0.007856892375950346,0.003838494918454322)
$data[33] = @(2,1,1.65374,3.30748,0.0459128248697883,0.03147194456382552,3,1,1.161330666666667,3.483992,0.02265142207290381,0.02421627468347122,1.920538976693333,11.52323386016,0.00103999077468489,0.0007621332543805775)
$data[34] = @(2,1,1.65374,3.30748,0.0459128248697883,0.03147194456382552,3,1,7.719567000000001,23.158701,0.1505679436150197,0.1609697911844745,12.76615673058,76.59694038348,0.006912999626200559,0.00506603234460835)
$data[35] = @(2,1,1.65374,3.30748,0.0459128248697883,0.03147194456382552,3,1,16.007757,48.023271,0.3122267160034931,0.3337966108230954,26.47266806118,158.83600836708,0.0143352105315375,0.0105052284314173)
$data[36] = @(2,1,1.65374,3.30748,0.0459128248697883,0.03147194456382552,3,1,16.441887,49.325661,0.3206942973278844,0.3428491671966479,27.19060620738,163.14363724428,0.01472398110995497,0.01079012998376665)
$data[37] = @(2,1,1.65374,3.30748,0.0459128248697883,0.03147194456382552,2,1,1.165529,2.331058,0.02273330936225701,0.01620254605380927,1.92748192846,7.70992771384,0.001043750451460025,0.0005099256311983151)

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}
